$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.237.82"
$ws.Range("E2").Value = "  +5.94%  "
$ws.Range("D3").Value = "2.286.48"
$ws.Range("E3").Value = "  +3.62%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.92"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.641"
$ws.Range("E6").Value = "  +3.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.19"
$ws.Range("E7").Value = "  +8.26%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.436"
$ws.Range("E9").Value = "  +8.56%  "
$ws.Range("E10").Value = "  +17.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.45"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.06"
$ws.Range("E12").Value = "  +17.73%  "
$ws.Range("E13").Value = "  +1.81%  "
$ws.Range("D14").Value = "2.624.90"
$ws.Range("E14").Value = "  +3.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.09"
$ws.Range("E15").Value = "  +4.79%  "
$ws.Range("E16").Value = "  +8.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.838"
$ws.Range("E17").Value = "  +5.90%  "
$ws.Range("D18").Value = "2.287.36"
$ws.Range("E18").Value = "  +3.67%  "
$ws.Range("D19").Value = "44.115.49"
$ws.Range("E19").Value = "  +5.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000101"
$ws.Range("E20").Value = "  +11.72%  "
$ws.Range("E21").Value = "  +3.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.12"
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "261.85"
$ws.Range("E23").Value = "  +8.44%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.18"
$ws.Range("E27").Value = "  +5.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.73"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "21.10"
$ws.Range("E29").Value = "  +6.91%  "
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.82"
$ws.Range("E32").Value = "  +7.49%  "
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0690"
$ws.Range("E34").Value = "  +6.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.76"
$ws.Range("E35").Value = "  +3.32%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  +8.76%  "
$ws.Range("E38").Value = "  +8.60%  "
$ws.Range("E39").Value = "  +1.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0249"
$ws.Range("E40").Value = "  +4.18%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.76"
$ws.Range("E43").Value = "  +9.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0973"
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.51"
$ws.Range("E45").Value = "  +1.62%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.46"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.21"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("D48").Value = "1.476.93"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.35"
$ws.Range("E49").Value = "  +6.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000210"
$ws.Range("E50").Value = "  -11.27%  "
$ws.Range("E51").Value = "  +2.18%  "
